$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("missing_values")
$ws1.Range("B4").Value = 180
$ws1.Range("C4").Value = 1.4410375470338643
$ws1.Range("B5").Value = 348
$ws1.Range("C5").Value = 2.7860059242654711
$ws1.Range("B6").Value = 1579
$ws1.Range("C6").Value = 12.641101593147066
$ws1.Range("B8").Value = 1927
$ws1.Range("C8").Value = 15.427107517412535
$ws1.Range("B9").Value = 10310
$ws1.Range("C9").Value = 82.539428388439674
$ws1.Range("B10").Value = 12491
$ws1.Range("B15").Value = 22499
$ws1.Range("C15").Value = 370.04934210526318
$ws1.Range("C16").Value = 0.1151315789473684
$ws1.Range("C17").Value = 0.2960526315789474
$ws1.Range("C19").Value = 0.41118421052631576
$ws1.Range("B20").Value = 5990
$ws1.Range("C20").Value = 98.51973684210526
$ws1.Range("B21").Value = 6080
$ws1.Range("C26").Value = 1.7201834862385321
$ws1.Range("B27").Value = 17
$ws1.Range("C27").Value = 1.9495412844036699
$ws1.Range("C28").Value = 0.22935779816513763
$ws1.Range("B29").Value = 19
$ws1.Range("C29").Value = 2.1788990825688073
$ws1.Range("B30").Value = 838
$ws1.Range("C30").Value = 96.100917431192656
$ws1.Range("B31").Value = 872
$ws1.Range("B36").Value = 24291
$ws1.Range("C36").Value = 315.9599375650364
$ws1.Range("B37").Value = 54
$ws1.Range("C37").Value = 0.70239334027055156
$ws1.Range("B39").Value = 54
$ws1.Range("C39").Value = 0.70239334027055156
$ws1.Range("B40").Value = 7627
$ws1.Range("C40").Value = 99.206555671175849
$ws1.Range("B41").Value = 7688

$ws2 = $wb.Worksheets.Item("profile_missing_values")
$ws2.Range("B4").Value = 31.551634665282823
$ws2.Range("B5").Value = 68.448365334717181
$ws2.Range("B7").Value = 0.36325895173845357
$ws2.Range("B8").Value = 12.454592631032693
$ws2.Range("B9").Value = 22.625843279709393
$ws2.Range("B10").Value = 23.196678775298391
$ws2.Range("B11").Value = 19.615983393876494
$ws2.Range("B12").Value = 14.737934613388687
$ws2.Range("B13").Value = 7.0057083549558907
$ws2.Range("B15").Value = 23.456149455111571
$ws2.Range("B16").Value = 33.004670472236633
$ws2.Range("B17").Value = 29.9429164504411
$ws2.Range("B18").Value = 10.949662688116243
$ws2.Range("B19").Value = 2.4909185262065385
$ws2.Range("B20").Value = 0.15568240788790866
$ws2.Range("B22").Value = 0.36325895173845357
$ws2.Range("B23").Value = 0.15568240788790866
$ws2.Range("B24").Value = 26.777374156720292
$ws2.Range("B25").Value = 48.209652309289055
$ws2.Range("B26").Value = 5.1894135962636225
$ws2.Range("B27").Value = 14.011416709911781
$ws2.Range("B28").Value = 0.31136481577581732
$ws2.Range("B29").Value = 4.9818370524130771
$ws2.Range("B31").Value = 4.4628956927867147
$ws2.Range("B32").Value = 4.3072132848988067
$ws2.Range("B33").Value = 11.572392319667877
$ws2.Range("B34").Value = 72.288531395952262
$ws2.Range("B35").Value = 0.36325895173845357
$ws2.Range("B36").Value = 4.0477426050856256
$ws2.Range("B37").Value = 2.7503892060197197
$ws2.Range("B38").Value = 0.20757654385054489
$ws2.Range("B40").Value = 18.266735858847948
$ws2.Range("B41").Value = 0.67462376751427089
$ws2.Range("B42").Value = 1.7644006227296314
$ws2.Range("B43").Value = 1.2454592631032693
$ws2.Range("B44").Value = 6.123508043591074
$ws2.Range("B45").Value = 20.550077841203944
$ws2.Range("B46").Value = 12.817851582771148
$ws2.Range("B47").Value = 2.9060716139076286
$ws2.Range("B48").Value = 6.0716139076284374
$ws2.Range("B49").Value = 29.372080954852102
$ws2.Range("B50").Value = 0.20757654385054489
$ws2.Range("B52").Value = 19.564089257913857
$ws2.Range("B53").Value = 80.435910742086151

$ws3 = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws3.Range("B4").Value = 2035627.05475047
$ws3.Range("C4").Value = 252179.71875
$ws3.Range("D4").Value = 453923.5
$ws3.Range("E4").Value = 1008718.875
$ws3.Range("G4").Value = 4034875.5
$ws3.Range("H4").Value = 2074505.4833767149
$ws3.Range("I4").Value = 252179.71875
$ws3.Range("J4").Value = 467037.0625
$ws3.Range("K4").Value = 1006662.8125
$ws3.Range("L4").Value = 2653803.75
$ws3.Range("M4").Value = 4102169.25

$ws4 = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws4.Range("B4").Value = 1527225.5785855504
$ws4.Range("D4").Value = 252179.71875
$ws4.Range("F4").Value = 2455983.5
$ws4.Range("G4").Value = 2556855.5
$ws4.Range("H4").Value = 1525197.7617169488
$ws4.Range("J4").Value = 252179.71875
$ws4.Range("L4").Value = 2455983.5
$ws4.Range("M4").Value = 2556855.5

$ws5 = $wb.Worksheets.Item("nonlabor_imp_stochastic_reg")
$ws5.Range("B4").Value = 9813130.2602175437
$ws5.Range("C4").Value = 121046.265625
$ws5.Range("D4").Value = 201743.78125
$ws5.Range("G4").Value = 942000
$ws5.Range("H4").Value = 9747674.9574410859
$ws5.Range("I4").Value = 120932.375
$ws5.Range("J4").Value = 201743.78125
$ws5.Range("M4").Value = 944749.9375

$ws6 = $wb.Worksheets.Item("labor_beneimp_stochastic_reg")
$ws6.Range("B4").Value = 873729.6433604639
$ws6.Range("C4").Value = 50435.9453125
$ws6.Range("D4").Value = 191656.59375
$ws6.Range("E4").Value = 302615.65625
$ws6.Range("F4").Value = 860816.1875
$ws6.Range("G4").Value = 2000000
$ws6.Range("H4").Value = 972859.7013830659
$ws6.Range("I4").Value = 55479.5390625
$ws6.Range("J4").Value = 196700.1875
$ws6.Range("K4").Value = 310968.96875
$ws6.Range("L4").Value = 870713.125
$ws6.Range("M4").Value = 2017437.75
